$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.08243498013860336
$ws.Range("J2").Value = 0.08243498013860337
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1311436666666667
$ws.Range("N2").Value = 0.393431
$ws.Range("O2").Value = 0.02663441993971509
$ws.Range("P2").Value = 0.02663441993971509
$ws.Range("Q2").Value = 0.01795780827855556
$ws.Range("R2").Value = 0.161620274507
$ws.Range("S2").Value = 0.002195607878733634
$ws.Range("T2").Value = 0.002195607878733635

$ws.Range("I3").Value = 0.08243498013860336
$ws.Range("J3").Value = 0.08243498013860337
$ws.Range("O3").Value = 0.06149297381279183
$ws.Range("P3").Value = 0.06149297381279183
$ws.Range("S3").Value = 0.005069172074921151
$ws.Range("T3").Value = 0.005069172074921152

$ws.Range("I4").Value = 0.08243498013860336
$ws.Range("J4").Value = 0.08243498013860337
$ws.Range("M4").Value = 4.009307333333333
$ws.Range("N4").Value = 12.027922
$ws.Range("O4").Value = 0.8142640654908683
$ws.Range("P4").Value = 0.8142640654908684
$ws.Range("Q4").Value = 0.5490038082037778
$ws.Range("R4").Value = 4.941034273834
$ws.Range("S4").Value = 0.06712384206631815
$ws.Range("T4").Value = 0.06712384206631816

$ws.Range("I5").Value = 0.08243498013860336
$ws.Range("J5").Value = 0.08243498013860337
$ws.Range("M5").Value = 0.480609
$ws.Range("N5").Value = 1.441827
$ws.Range("O5").Value = 0.09760854075662465
$ws.Range("P5").Value = 0.09760854075662465
$ws.Range("Q5").Value = 0.06581091179100002
$ws.Range("R5").Value = 0.592298206119
$ws.Range("S5").Value = 0.008046358118630409
$ws.Range("T5").Value = 0.008046358118630411

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.524162666666667
$ws.Range("H6").Value = 4.572488
$ws.Range("I6").Value = 0.9175650198613966
$ws.Range("J6").Value = 0.9175650198613967
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1311436666666667
$ws.Range("N6").Value = 0.393431
$ws.Range("O6").Value = 0.02663441993971509
$ws.Range("P6").Value = 0.02663441993971509
$ws.Range("Q6").Value = 0.1998842807031111
$ws.Range("R6").Value = 1.798958526328
$ws.Range("S6").Value = 0.02443881206098145
$ws.Range("T6").Value = 0.02443881206098146

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.524162666666667
$ws.Range("H7").Value = 4.572488
$ws.Range("I7").Value = 0.9175650198613966
$ws.Range("J7").Value = 0.9175650198613967
$ws.Range("O7").Value = 0.06149297381279183
$ws.Range("P7").Value = 0.06149297381279183
$ws.Range("Q7").Value = 0.4614885124844444
$ws.Range("R7").Value = 4.15339661236
$ws.Range("S7").Value = 0.05642380173787068
$ws.Range("T7").Value = 0.05642380173787068

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.524162666666667
$ws.Range("H8").Value = 4.572488
$ws.Range("I8").Value = 0.9175650198613966
$ws.Range("J8").Value = 0.9175650198613967
$ws.Range("M8").Value = 4.009307333333333
$ws.Range("N8").Value = 12.027922
$ws.Range("O8").Value = 0.8142640654908683
$ws.Range("P8").Value = 0.8142640654908684
$ws.Range("Q8").Value = 6.110836556659555
$ws.Range("R8").Value = 54.997529009936
$ws.Range("S8").Value = 0.7471402234245501
$ws.Range("T8").Value = 0.7471402234245503

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.524162666666667
$ws.Range("H9").Value = 4.572488
$ws.Range("I9").Value = 0.9175650198613966
$ws.Range("J9").Value = 0.9175650198613967
$ws.Range("M9").Value = 0.480609
$ws.Range("N9").Value = 1.441827
$ws.Range("O9").Value = 0.09760854075662465
$ws.Range("P9").Value = 0.09760854075662465
$ws.Range("Q9").Value = 0.7325262950639999
$ws.Range("R9").Value = 6.592736655576
$ws.Range("S9").Value = 0.08956218263799423
$ws.Range("T9").Value = 0.08956218263799424
